# Adds a new "percentage_change" column (AB) to the sheet, and fills in
# the target_wval / reference_wval (Z/AA) "corner" values that were
# previously missing for a handful of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header for column AB (match the header style used by the
#    rest of row 1, e.g. bold / centered / bordered).
# ---------------------------------------------------------------------
$ws.Range("AA1").Copy() | Out-Null
$ws.Range("AB1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("AB1").Value = "percentage_change"

# ---------------------------------------------------------------------
# 2. target_wval (Z) / reference_wval (AA) values that need to be filled
#    in for the rows that previously had no "corner" data.
# ---------------------------------------------------------------------
$zaa = @{
    9  = @{ Z = 0.290291;            AA = 0.296605 }
    10 = @{ Z = 3.444819164217974;   AA = 3.371487331636351 }
    23 = @{ Z = 0.08673139999999999; AA = 0.09739480000000002 }
    24 = @{ Z = 11.52984962770116;   AA = 10.26748861335513 }
    30 = @{ Z = 0.259733;            AA = 0.283546 }
    31 = @{ Z = 3.850107610507713;   AA = 3.526764616675953 }
}

foreach ($r in $zaa.Keys) {
    $ws.Range("Z$r").Value = $zaa[$r].Z
    $ws.Range("AA$r").Value = $zaa[$r].AA
}

# ---------------------------------------------------------------------
# 3. percentage_change (AB) values -- only the "corner" rows (the first
#    two rows of each 7-row IP block) get a computed value; every other
#    row just gets a blank, formatted cell (matching column Z/AA).
# ---------------------------------------------------------------------
$ab = @{
    2  = -5.452515766698875
    3  = 5.766960179759497
    9  = -2.128757101195188
    10 = 2.175058820287223
    16 = -10.13073792006415
    17 = 11.27275075548432
    23 = -10.94863380796514
    24 = 12.2947398520029
    30 = -8.398284581690438
    31 = 9.168261252902004
}

for ($r = 2; $r -le 36; $r++) {
    # Match the formatting already used by the neighbouring Z column
    # (which carries the correct per-row fill style) before writing
    # the value (or leaving the cell blank).
    $ws.Range("Z$r").Copy() | Out-Null
    $ws.Range("AB$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    if ($ab.ContainsKey($r)) {
        $ws.Range("AB$r").Value = $ab[$r]
    }
}
